$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 305, pushing existing rows 305-326 down to 306-327
$ws.Rows.Item(305).Insert()

# Populate the newly inserted row 305 with the new record's data.
# Non-numeric/date columns mirror the values used by the surrounding rows.
$ws.Range("A305").Value = 4
$ws.Range("B305").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C305").Value = "Los Lagos"
$ws.Range("D305").Value = 44826
$ws.Range("E305").Value = 10
$ws.Range("F305").Value = 100112043
$ws.Range("G305").Value = "Pepino ensalada"
$ws.Range("H305").Value = "Sin especificar"
$ws.Range("I305").Value = "Primera"
$ws.Range("J305").Value = 200
$ws.Range("K305").Value = 26000
$ws.Range("L305").Value = 26000
$ws.Range("M305").Value = 26000
$ws.Range("N305").Value = "$/caja 60 unidades"
$ws.Range("O305").Value = "Región de Arica y Parinacota"
$ws.Range("P305").Value = 433
$ws.Range("Q305").Value = 60
$ws.Range("R305").Value = "Hortaliza"
